$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.05636357315966059
$ws.Range("C2").Value = 0.9992368720115482
$ws.Range("D2").Value = 0.1914382525850619
$ws.Range("F2").Value = "Pipeline(steps=[('model',`n                 RandomForestRegressor(max_depth=5, n_estimators=150))])"
$ws.Range("G2").Value = 0.124317388383497
$ws.Range("H2").Value = 0.991
